$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per the meeting, the "Meetings" entry in the Nov 20th - Nov 27th timesheet
# block dropped from a full day (1) down to 2 hours -- expressed as a time
# fraction of a day (8.3333...E-2), matching the h:mm formatting already used
# by the other time cells in that block (e.g. C25/C26).
$ws.Range("C27").Value = 0.083333333333333329
$ws.Range("C27").NumberFormat = $ws.Range("C26").NumberFormat

# Reflect where the user was working when they made the change: the weekly
# total rows near the bottom of the sheet.
$ws.Range("E31:E32").Select()
